# Revert "feat(dialog): update CN data and dialogue Excel files"
# This removes the rows that were added by the original feature commit:
#   - row with id "kettle_eye"
#   - rows with ids "sorin", "issizzle", "vishnu", "az"
# Deleting the rows causes Excel to drop the now-unreferenced shared
# strings automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Person")

# Delete from the bottom up so row numbers of earlier rows are not
# shifted out from under us while we still need them.

# Rows 37-40: sorin, issizzle, vishnu, az
$ws.Rows.Item(37).Resize(4).EntireRow.Delete()

# Row 8: kettle_eye
$ws.Rows.Item(8).EntireRow.Delete()
